$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.218.90"
$ws.Range("E2").Value = "  -1.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.600.59"
$ws.Range("E3").Value = "  -0.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.44"
$ws.Range("E5").Value = "  +2.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.75"
$ws.Range("E6").Value = "  -0.39%  "

# Row 7
$ws.Range("E7").Value = "  +0.32%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.54"
$ws.Range("E9").Value = "  -2.33%  "

# Row 10
$ws.Range("E10").Value = "  -1.30%  "

# Row 11
$ws.Range("E11").Value = "  -1.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.370"
$ws.Range("E12").Value = "  -1.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.061.20"
$ws.Range("E13").Value = "  -0.55%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.49"
$ws.Range("E14").Value = "  +3.85%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.220.97"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.603.97"
$ws.Range("E17").Value = "  -0.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.33"
$ws.Range("E18").Value = "  +2.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.61"
$ws.Range("E19").Value = "  -1.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.58"
$ws.Range("E20").Value = "  -1.14%  "

# Row 21
$ws.Range("E21").Value = "  -2.54%  "

# Row 22
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("E23").Value = "  +2.80%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.75"
$ws.Range("E24").Value = "  +0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.15%  "

# Row 26
$ws.Range("E26").Value = "  -1.55%  "

# Row 27
$ws.Range("E27").Value = "  +3.01%  "

# Row 28
$ws.Range("E28").Value = "  +2.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30: 'Aptos' -> 'Monero'
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.58"
$ws.Range("E30").Value = "  +4.35%  "

# Row 31: 'USDe' -> 'Aptos'
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.40"
$ws.Range("E31").Value = "  +1.92%  "

# Row 32: 'Monero' -> 'USDe'
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +0.13%  "

# Row 33
$ws.Range("E33").Value = "  -0.41%  "

# Row 34
$ws.Range("E34").Value = "  +8.75%  "

# Row 35
$ws.Range("E35").Value = "  +0.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.981"
$ws.Range("E36").Value = "  +1.62%  "

# Row 37
$ws.Range("E37").Value = "  +2.71%  "

# Row 38
$ws.Range("E38").Value = "  +1.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "314.29"
$ws.Range("E39").Value = "  +3.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.87"
$ws.Range("E40").Value = "  +1.71%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.841"
$ws.Range("E41").Value = "  -2.03%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.61"
$ws.Range("E42").Value = "  -3.90%  "

# Row 43
$ws.Range("E43").Value = "  +0.81%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.62%  "

# Row 45
$ws.Range("E45").Value = "  +1.14%  "

# Row 46
$ws.Range("E46").Value = "  -0.30%  "

# Row 47
$ws.Range("E47").Value = "  -0.02%  "

# Row 48: 'RenderToken' -> 'VeChain'
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0242"
$ws.Range("E48").Value = "  -0.14%  "

# Row 49: 'VeChain' -> 'RenderToken'
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.96"
$ws.Range("E49").Value = "  +2.70%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.86"
$ws.Range("E50").Value = "  +1.68%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.74"
$ws.Range("E51").Value = "  +0.41%  "
